$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header title text in B1 (drop the period after "6.4.2.1")
$ws.Range("B1").Value = "6.4.2.1 Общий объем забора пресной воды "

# Update data values for the 2022 column (L)
$ws.Range("L5").Value = 8741.9
$ws.Range("L7").Value = 8483.5
$ws.Range("L14").Value = 1327.6
$ws.Range("L18").Value = 54

# Move the active selection to O2
$ws.Range("O2").Select()
